$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 74.86362832748885
$ws.Range("B3").Value = 77.26092484587997
$ws.Range("B4").Value = 84.3796513590588
$ws.Range("H5").Value = 96.3909448788188
$ws.Range("H6").Value = 96.16231408458287
$ws.Range("H7").Value = 96.22258537269308
$ws.Range("C8").Value = 96.40155686798889
$ws.Range("C9").Value = 94.44352071792737
$ws.Range("C10").Value = 95.49943615400775
$ws.Range("D11").Value = 98.87493468609266
$ws.Range("D12").Value = 98.95576412193684
$ws.Range("D13").Value = 98.99581633356149
$ws.Range("E14").Value = 98.67938263787825
$ws.Range("E15").Value = 98.6125584401299
$ws.Range("E16").Value = 98.65159386129629
$ws.Range("F17").Value = 98.22722310680443
$ws.Range("F18").Value = 98.31651001515212
$ws.Range("F19").Value = 98.3234674102959
$ws.Range("G20").Value = 97.49108089310863
$ws.Range("G21").Value = 97.64167167066532
$ws.Range("G22").Value = 97.61444244323198
$ws.Range("B23").Value = 83.80384211316442
$ws.Range("B24").Value = 85.81581984353001
$ws.Range("H25").Value = 96.42074442171868
$ws.Range("H26").Value = 96.05109446553321
$ws.Range("C27").Value = 94.71593242995539
$ws.Range("C28").Value = 94.42800757941872
$ws.Range("D29").Value = 99.1950229664282
$ws.Range("D30").Value = 99.04699111772418
$ws.Range("E31").Value = 98.71822783560339
$ws.Range("E32").Value = 98.57715618753579
$ws.Range("F33").Value = 98.21039610442436
$ws.Range("F34").Value = 98.3988521742662
$ws.Range("G35").Value = 97.56282852942
$ws.Range("G36").Value = 97.49849673505496
$ws.Range("B37").Value = 83.07767861872912
$ws.Range("B38").Value = 90.86934023608504
$ws.Range("H39").Value = 96.20008212163397
$ws.Range("H40").Value = 96.19117178318243
$ws.Range("C41").Value = 96.02638804220749
$ws.Range("C42").Value = 95.97763622075473
$ws.Range("D43").Value = 99.16235665393695
$ws.Range("D44").Value = 98.9887433767552
$ws.Range("E45").Value = 98.68142911248921
$ws.Range("E46").Value = 98.4069337158529
$ws.Range("F47").Value = 98.29531309204827
$ws.Range("F48").Value = 98.32922777190208
$ws.Range("G49").Value = 97.5079218435293
$ws.Range("G50").Value = 97.60429918917126
